$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "age"
$ws.Range("A2").Value = "suraj"
$ws.Range("B2").Value = 30
$ws.Range("A3").Value = "rohan"
$ws.Range("B3").Value = 32
$ws.Range("A4").Value = "ashish"
$ws.Range("B4").Value = 34

$ws.Range("B5").Select()
